{"js": "// Apply the \"fruity-gold\" rewrite: update the title (appears twice),\n// the four \"what we like\" bullets, the two \"what we don't like\" bullets,\n// and the closing bold title / italic summary lines.\n\nconst replacements = [\n  [\n    \"Play Fruity Gold Free: Review of the Traditional Video Slot Game\",\n    \"Play Fruity Gold for Free - Traditional Slot Game\",\n  ],\n  [\n    \"Traditional and straightforward gameplay\",\n    \"Straightforward gameplay with no surprises or special rounds\",\n  ],\n  [\n    \"Double Wilds and Wilds offer ample opportunities to win\",\n    \"Ample opportunities to win significant prizes with multipliers and Wilds\",\n  ],\n  [\n    \"Nice graphics and familiar fruit symbols\",\n    \"Traditional fruit symbols that many players love and appreciate\",\n  ],\n  [\n    \"27 paylines offer more chances to win\",\n    \"Brings back the old typical atmosphere of a real casino\",\n  ],\n  [\n    \"No extra or free spins\",\n    \"Lack of extra spins, free spins, animations, or any special functions\",\n  ],\n  [\n    \"Graphics are a bit too minimal\",\n    \"Graphics and sounds may be too minimal and lacking in detail\",\n  ],\n  [\n    \"Discover the traditional gameplay and fruit symbols in Fruity Gold. Play the game for free and read our review for information on features and graphics.\",\n    \"Experience the traditional atmosphere of a real casino in Fruity Gold. Play now for free!\",\n  ],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const result of results.items) {\n    result.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the \"fruity-gold\" rewrite: update the title (appears twice),\n# the four \"what we like\" bullets, the two \"what we don't like\" bullets,\n# and the closing bold title / italic summary lines.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"Play Fruity Gold Free: Review of the Traditional Video Slot Game\", \"Play Fruity Gold for Free - Traditional Slot Game\"),\n    @(\"Traditional and straightforward gameplay\", \"Straightforward gameplay with no surprises or special rounds\"),\n    @(\"Double Wilds and Wilds offer ample opportunities to win\", \"Ample opportunities to win significant prizes with multipliers and Wilds\"),\n    @(\"Nice graphics and familiar fruit symbols\", \"Traditional fruit symbols that many players love and appreciate\"),\n    @(\"27 paylines offer more chances to win\", \"Brings back the old typical atmosphere of a real casino\"),\n    @(\"No extra or free spins\", \"Lack of extra spins, free spins, animations, or any special functions\"),\n    @(\"Graphics are a bit too minimal\", \"Graphics and sounds may be too minimal and lacking in detail\"),\n    @(\"Discover the traditional gameplay and fruit symbols in Fruity Gold. Play the game for free and read our review for information on features and graphics.\", \"Experience the traditional atmosphere of a real casino in Fruity Gold. Play now for free!\")\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n    $d.Content.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n"}
